# Auto-applied scheduled-runner price/profit refresh for Typhon_Profits leve sheets.
# Updates currentAveragePrice* / Leve*Price* / LeveProfit* columns (H:N) per row,
# mirroring an external market-data pull. Values/removed cells match the upstream diff.

$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 5 (Leve Item ID 5503)
$ws.Range("H5").Value = 909.8333
$ws.Range("I5").Value = 14.75
$ws.Range("J5").Value = 2700
$ws.Range("K5").Value = 14.75
$ws.Range("L5").Value = 2700
$ws.Range("M5").Value = 100.25
$ws.Range("N5").Value = -2930
# Row 32 (Leve Item ID 5484)
$ws.Range("H32").Value = 272.75
$ws.Range("J32").Value = 272.75
$ws.Range("L32").Value = 272.75
$ws.Range("N32").Value = -924.75
# Row 33 (Leve Item ID 5512)
$ws.Range("H33").Value = 153.71428
$ws.Range("I33").Value = 163.23077
$ws.Range("K33").Value = 163.23077
$ws.Range("M33").Value = 65.76922999999999
# Row 99 (Leve Item ID 19883)
$ws.Range("H99").Value = 202.33333
$ws.Range("I99").Value = 192.8
$ws.Range("K99").Value = 578.4000000000001
$ws.Range("M99").Value = 919.5999999999999
# Row 129 (Leve Item ID 36115)
$ws.Range("H129").Value = 263891.1
$ws.Range("J129").Value = 294896.25
$ws.Range("L129").Value = 884688.75
$ws.Range("N129").Value = -894688.75
# Row 132 (Leve Item ID 44049)
$ws.Range("H132").Value = 3080.125
$ws.Range("I132").Value = 3391.2593
$ws.Range("K132").Value = 10173.7779
$ws.Range("M132").Value = -7643.777900000001
# Row 138 (Leve Item ID 44169)
$ws.Range("H138").Value = 2139.0312
$ws.Range("J138").Value = 2595.5952
$ws.Range("L138").Value = 7786.785600000001
$ws.Range("N138").Value = -18066.7856

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32 (Leve Item ID 44147)
$ws.Range("H32").Value = 5530.153
$ws.Range("I32").Value = 4283.62
$ws.Range("K32").Value = 4283.62
$ws.Range("M32").Value = -3996.62
# Row 40 (Leve Item ID 3833)
$ws.Range("H40").Value = 0
$ws.Range("J40").Value = 0
$ws.Range("L40").Value = 0
$ws.Range("N40").ClearContents()

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 54 (Leve Item ID 2376)
$ws.Range("H54").Value = 8659.333000000001
$ws.Range("I54").Value = 2987.4285
$ws.Range("J54").Value = 16600
$ws.Range("K54").Value = 2987.4285
$ws.Range("L54").Value = 16600
$ws.Range("M54").Value = -2503.4285
$ws.Range("N54").Value = -17568
# Row 95 (Leve Item ID 18194)
$ws.Range("H95").Value = 19541.334
$ws.Range("J95").Value = 19541.334
$ws.Range("L95").Value = 19541.334
$ws.Range("N95").Value = -25033.334
# Row 107 (Leve Item ID 27706)
$ws.Range("H107").Value = 2032.6666
$ws.Range("J107").Value = 2100
$ws.Range("L107").Value = 2100
$ws.Range("N107").Value = -5940

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 8 (Leve Item ID 1894)
$ws.Range("H8").Value = 6312
$ws.Range("J8").Value = 7140
$ws.Range("L8").Value = 7140
$ws.Range("N8").Value = -7420
# Row 31 (Leve Item ID 44023)
$ws.Range("H31").Value = 4346.241
$ws.Range("I31").Value = 2309.125
$ws.Range("J31").Value = 6853.4614
$ws.Range("K31").Value = 2309.125
$ws.Range("L31").Value = 6853.4614
$ws.Range("M31").Value = -2014.125
$ws.Range("N31").Value = -7443.4614
# Row 34 (Leve Item ID 44023)
$ws.Range("H34").Value = 4346.241
$ws.Range("I34").Value = 2309.125
$ws.Range("J34").Value = 6853.4614
$ws.Range("K34").Value = 2309.125
$ws.Range("L34").Value = 6853.4614
$ws.Range("M34").Value = -2107.125
$ws.Range("N34").Value = -7257.4614
# Row 86 (Leve Item ID 12584)
$ws.Range("H86").Value = 24885
$ws.Range("J86").Value = 29813.2
$ws.Range("L86").Value = 29813.2
$ws.Range("N86").Value = -32059.2
# Row 89 (Leve Item ID 12584)
$ws.Range("H89").Value = 24885
$ws.Range("J89").Value = 29813.2
$ws.Range("L89").Value = 149066
$ws.Range("N89").Value = -160298
# Row 99 (Leve Item ID 36198)
$ws.Range("H99").Value = 3261.56
$ws.Range("I99").Value = 2761.1765
$ws.Range("J99").Value = 4324.875
$ws.Range("K99").Value = 2761.1765
$ws.Range("L99").Value = 4324.875
$ws.Range("M99").Value = -1263.1765
$ws.Range("N99").Value = -7320.875
# Row 126 (Leve Item ID 36198)
$ws.Range("H126").Value = 3261.56
$ws.Range("I126").Value = 2761.1765
$ws.Range("J126").Value = 4324.875
$ws.Range("K126").Value = 8283.529500000001
$ws.Range("L126").Value = 12974.625
$ws.Range("M126").Value = -5813.529500000001
$ws.Range("N126").Value = -17914.625

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 97 (Leve Item ID 19846)
$ws.Range("H97").Value = 614.75
$ws.Range("J97").Value = 614.75
$ws.Range("L97").Value = 1844.25
$ws.Range("N97").Value = -2836.25
# Row 98 (Leve Item ID 19843)
$ws.Range("H98").Value = 863.3333
$ws.Range("I98").Value = 398.75
$ws.Range("J98").Value = 1235
$ws.Range("K98").Value = 1196.25
$ws.Range("L98").Value = 3705
$ws.Range("M98").Value = 301.75
$ws.Range("N98").Value = -6701
# Row 117 (Leve Item ID 27870)
$ws.Range("H117").Value = 1007.86664
$ws.Range("I117").Value = 847.4
$ws.Range("J117").Value = 1088.1
$ws.Range("K117").Value = 2542.2
$ws.Range("L117").Value = 3264.3
$ws.Range("M117").Value = 899.8000000000002
$ws.Range("N117").Value = -10148.3
# Row 131 (Leve Item ID 36060)
$ws.Range("H131").Value = 718.84
$ws.Range("J131").Value = 741.44086
$ws.Range("L131").Value = 2224.32258
$ws.Range("N131").Value = -12304.32258
# Row 140 (Leve Item ID 44097)
$ws.Range("H140").Value = 1859.2941
$ws.Range("I140").Value = 911.5714
$ws.Range("J140").Value = 3390.2307
$ws.Range("K140").Value = 2734.7142
$ws.Range("L140").Value = 10170.6921
$ws.Range("M140").Value = 2445.2858
$ws.Range("N140").Value = -20530.6921

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 46 (Leve Item ID 2078)
$ws.Range("H46").Value = 20364.143
$ws.Range("J46").Value = 20424.834
$ws.Range("L46").Value = 20424.834
$ws.Range("N46").Value = -20736.834
# Row 69 (Leve Item ID 11891)
$ws.Range("H69").Value = 0
$ws.Range("J69").Value = 0
$ws.Range("L69").Value = 0
$ws.Range("N69").ClearContents()
# Row 70 (Leve Item ID 14146)
$ws.Range("H70").Value = 3916244.2
$ws.Range("I70").Value = 4488.5
$ws.Range("K70").Value = 4488.5
$ws.Range("M70").Value = -4218.5
# Row 72 (Leve Item ID 11891)
$ws.Range("H72").Value = 0
$ws.Range("J72").Value = 0
$ws.Range("L72").Value = 0
$ws.Range("N72").ClearContents()
# Row 73 (Leve Item ID 14146)
$ws.Range("H73").Value = 3916244.2
$ws.Range("I73").Value = 4488.5
$ws.Range("K73").Value = 4488.5
$ws.Range("M73").Value = -3552.5

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 7 (Leve Item ID 36249)
$ws.Range("H7").Value = 3545
$ws.Range("I7").Value = 3697.5
$ws.Range("J7").Value = 3179
$ws.Range("K7").Value = 3697.5
$ws.Range("L7").Value = 3179
$ws.Range("M7").Value = -3585.5
$ws.Range("N7").Value = -3403
# Row 24 (Leve Item ID 3774)
$ws.Range("H24").Value = 2577.8
$ws.Range("J24").Value = 2577.8
$ws.Range("L24").Value = 2577.8
$ws.Range("N24").Value = -3263.8
# Row 25 (Leve Item ID 3547)
$ws.Range("H25").Value = 4004
$ws.Range("J25").Value = 4004
$ws.Range("L25").Value = 4004
$ws.Range("N25").Value = -4464
# Row 55 (Leve Item ID 5284)
$ws.Range("H55").Value = 778.46155
$ws.Range("I55").Value = 1184.75
$ws.Range("J55").Value = 128.4
$ws.Range("K55").Value = 1184.75
$ws.Range("L55").Value = 128.4
$ws.Range("M55").Value = -1011.75
$ws.Range("N55").Value = -474.4
# Row 126 (Leve Item ID 36249)
$ws.Range("H126").Value = 3545
$ws.Range("I126").Value = 3697.5
$ws.Range("J126").Value = 3179
$ws.Range("K126").Value = 11092.5
$ws.Range("L126").Value = 9537
$ws.Range("M126").Value = -8622.5
$ws.Range("N126").Value = -14477

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 76 (Leve Item ID 10896)
$ws.Range("H76").Value = 37000
$ws.Range("J76").Value = 37000
$ws.Range("L76").Value = 37000
$ws.Range("N76").Value = -37630
# Row 79 (Leve Item ID 10896)
$ws.Range("H79").Value = 37000
$ws.Range("J79").Value = 37000
$ws.Range("L79").Value = 37000
$ws.Range("N79").Value = -39184
# Row 86 (Leve Item ID 11977)
$ws.Range("H86").Value = 22000
$ws.Range("J86").Value = 22000
$ws.Range("L86").Value = 22000
$ws.Range("N86").Value = -24246
# Row 89 (Leve Item ID 11977)
$ws.Range("H89").Value = 22000
$ws.Range("J89").Value = 22000
$ws.Range("L89").Value = 110000
$ws.Range("N89").Value = -121232
# Row 92 (Leve Item ID 18088)
$ws.Range("H92").Value = 23749.5
$ws.Range("J92").Value = 23749.5
$ws.Range("L92").Value = 23749.5
$ws.Range("N92").Value = -28741.5
# Row 95 (Leve Item ID 18243)
$ws.Range("H95").Value = 27699.8
$ws.Range("J95").Value = 27699.8
$ws.Range("L95").Value = 27699.8
$ws.Range("N95").Value = -33191.8
